# Add a new column N (year 2022) to the table, mirroring the formatting
# already used for column M, and update the sheet's selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: blank, bottom-bordered cell (same style as M2), no value.
$ws.Range("M2").Copy($ws.Range("N2"))

# Row 3: year header value, same style as M3 (no special number format).
$ws.Range("M3").Copy($ws.Range("N3"))
$ws.Range("N3").Value = 2022

# Row 4: numeric value, style already carries the "0.0" number format
# (same as M4), so just copy style + set value.
$ws.Range("M4").Copy($ws.Range("N4"))
$ws.Range("N4").Value = 9.224468514531754

# Rows 5-12: numeric values whose style is the same as the matching M-column
# style but with the "0.0" number format applied.
$ws.Range("M5").Copy($ws.Range("N5"))
$ws.Range("N5").Value = 4.6068543125097872
$ws.Range("N5").NumberFormat = "0.0"

$ws.Range("M6").Copy($ws.Range("N6"))
$ws.Range("N6").Value = 13.543910285971602
$ws.Range("N6").NumberFormat = "0.0"

$ws.Range("M7").Copy($ws.Range("N7"))
$ws.Range("N7").Value = 24.703327617190443
$ws.Range("N7").NumberFormat = "0.0"

$ws.Range("M8").Copy($ws.Range("N8"))
$ws.Range("N8").Value = 28.608474183838851
$ws.Range("N8").NumberFormat = "0.0"

$ws.Range("M9").Copy($ws.Range("N9"))
$ws.Range("N9").Value = 20.904451081350146
$ws.Range("N9").NumberFormat = "0.0"

$ws.Range("M10").Copy($ws.Range("N10"))
$ws.Range("N10").Value = 26.720095429750884
$ws.Range("N10").NumberFormat = "0.0"

$ws.Range("M11").Copy($ws.Range("N11"))
$ws.Range("N11").Value = 27.704327204727914
$ws.Range("N11").NumberFormat = "0.0"

$ws.Range("M12").Copy($ws.Range("N12"))
$ws.Range("N12").Value = 25.731792255708452
$ws.Range("N12").NumberFormat = "0.0"

# Update the selected cell to match the target workbook.
$ws.Range("Q5").Select()
